$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1223-1224, shifting the existing rows
# (old 1223 onward) down by two. This matches the change in <dimension>
# from A1:R1337 to A1:R1339.
$ws.Range("A1223:R1224").EntireRow.Insert()

# New row 1223: Primera, dated 2023-08-28 (serial 45166)
$ws.Range("A1223").Value = 6
$ws.Range("B1223").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1223").Value = "Metropolitana"
$ws.Range("D1223").Value = 45166
$ws.Range("E1223").Value = 13
$ws.Range("F1223").Value = 100112009
$ws.Range("G1223").Value = "Acelga"
$ws.Range("H1223").Value = "Sin especificar"
$ws.Range("I1223").Value = "Primera"
$ws.Range("J1223").Value = 190
$ws.Range("K1223").Value = 12000
$ws.Range("L1223").Value = 12000
$ws.Range("M1223").Value = 12000
$ws.Range("N1223").Value = "`$/docena de atados"
$ws.Range("O1223").Value = "Región Metropolitana"
$ws.Range("P1223").Value = 4000
$ws.Range("Q1223").Value = 3
$ws.Range("R1223").Value = "Hortaliza"

# New row 1224: Segunda, same date as row 1223
$ws.Range("A1224").Value = 6
$ws.Range("B1224").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1224").Value = "Metropolitana"
$ws.Range("D1224").Value = 45166
$ws.Range("E1224").Value = 13
$ws.Range("F1224").Value = 100112009
$ws.Range("G1224").Value = "Acelga"
$ws.Range("H1224").Value = "Sin especificar"
$ws.Range("I1224").Value = "Segunda"
$ws.Range("J1224").Value = 120
$ws.Range("K1224").Value = 9000
$ws.Range("L1224").Value = 9000
$ws.Range("M1224").Value = 9000
$ws.Range("N1224").Value = "`$/docena de atados"
$ws.Range("O1224").Value = "Región Metropolitana"
$ws.Range("P1224").Value = 3000
$ws.Range("Q1224").Value = 3
$ws.Range("R1224").Value = "Hortaliza"
